$wb = $excel.ActiveWorkbook

# Add the new "data laptop" worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "data laptop"

# Populate the values
$ws.Range("A1").Value = "datalaptop"
$ws.Range("A2").Value = "LAPTOP A"
$ws.Range("A3").Value = "LAPTOP B"
$ws.Range("A4").Value = "KOMPUTER"

# Column width to match original formatting (closest achievable quantized value)
$ws.Columns.Item(1).ColumnWidth = 20.65

# Selection state on the new sheet
$ws.Range("D7").Select()

$wb.Save()
